$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.781.07"
$ws.Range("E2").Value = "'  +0.46%  "
$ws.Range("D3").Value = "'1.649.71"
$ws.Range("E3").Value = "'  +1.05%  "
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("D5").Value = "'214.82"
$ws.Range("E5").Value = "'  +0.77%  "
$ws.Range("E6").Value = "'  +2.11%  "
$ws.Range("E7").Value = "'  -0.06%  "
$ws.Range("E8").Value = "'  -0.32%  "
$ws.Range("E9").Value = "'  +0.48%  "
$ws.Range("D10").Value = "'19.20"
$ws.Range("E10").Value = "'  +1.00%  "
$ws.Range("E11").Value = "'  +0.04%  "
$ws.Range("D12").Value = "'1.874.95"
$ws.Range("E12").Value = "'  +0.79%  "
$ws.Range("D13").Value = "'1.659.87"
$ws.Range("E13").Value = "'  +2.07%  "
$ws.Range("E14").Value = "'  +1.96%  "
$ws.Range("D15").Value = "'0.532"
$ws.Range("E15").Value = "'  +1.18%  "
$ws.Range("D16").Value = "'65.91"
$ws.Range("E16").Value = "'  +4.31%  "
$ws.Range("D17").Value = "'26.766.14"
$ws.Range("D18").Value = "'0.0₃0747"
$ws.Range("E18").Value = "'  +0.87%  "
$ws.Range("D19").Value = "'218.86"
$ws.Range("E19").Value = "'  +4.07%  "
$ws.Range("E20").Value = "'  -0.06%  "
$ws.Range("E21").Value = "'  +1.28%  "
$ws.Range("D22").Value = "'6.35"
$ws.Range("E22").Value = "'  +2.35%  "
$ws.Range("D23").Value = "'9.47"
$ws.Range("E23").Value = "'  +0.26%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "'  +10.43%  "
$ws.Range("D25").Value = "'147.75"
$ws.Range("E25").Value = "'  +0.42%  "
$ws.Range("E26").Value = "'  -0.03%  "
$ws.Range("E27").Value = "'  -0.03%  "
$ws.Range("D28").Value = "'6.95"
$ws.Range("E28").Value = "'  +0.81%  "
$ws.Range("D29").Value = "'15.87"
$ws.Range("E29").Value = "'  +3.09%  "
$ws.Range("E30").Value = "'  -0.46%  "
$ws.Range("E31").Value = "'  -0.31%  "
$ws.Range("D32").Value = "'3.39"
$ws.Range("E32").Value = "'  +4.45%  "
$ws.Range("E33").Value = "'  +2.89%  "
$ws.Range("D34").Value = "'1.272.02"
$ws.Range("E34").Value = "'  +8.68%  "
$ws.Range("E35").Value = "'  +1.36%  "
$ws.Range("D36").Value = "'2.38"
$ws.Range("E36").Value = "'  +0.93%  "
$ws.Range("D37").Value = "'0.0177"
$ws.Range("E37").Value = "'  +3.41%  "
$ws.Range("D38").Value = "'0.811"
$ws.Range("E38").Value = "'  -0.05%  "
$ws.Range("E39").Value = "'  +1.84%  "
$ws.Range("E40").Value = "'  -0.08%  "
$ws.Range("D41").Value = "'2.29"
$ws.Range("E41").Value = "'  -1.35%  "
$ws.Range("D42").Value = "'0.807"
$ws.Range("E42").Value = "'  +1.56%  "
$ws.Range("E43").Value = "'  -0.20%  "
$ws.Range("D44").Value = "'1.785.04"
$ws.Range("E44").Value = "'  +0.73%  "
$ws.Range("D45").Value = "'93.84"
$ws.Range("E45").Value = "'  +1.28%  "
$ws.Range("E46").Value = "'  +3.74%  "
$ws.Range("D47").Value = "'55.70"
$ws.Range("E47").Value = "'  +1.97%  "
$ws.Range("D48").Value = "'0.0₆0102"
$ws.Range("E48").Value = "'  -2.65%  "
$ws.Range("D50").Value = "'7.68"
$ws.Range("E50").Value = "'  +1.25%  "
$ws.Range("D51").Value = "'0.0969"
$ws.Range("E51").Value = "'  +3.11%  "

Write-Host "Updated 78 cells (D/E columns) for crypto price refresh."
